$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.821.44"
$ws.Range("D3").Value = "1.780.28"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.05"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0677"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "2.036.25"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("D14").Value = "1.787.01"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "33.850.42"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.58"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0510"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "1.382.86"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.644"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.60"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +13.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.66"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0139"
$ws.Range("E45").Value = "  +13.90%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0507"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.77"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").Value = "1.937.80"
$ws.Range("E50").Value = "  -0.89%  "
